$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H
$ws.Range("H1").EntireColumn.Insert()

# Set header for new column
$ws.Range("H1").Value = "interval"

# Set value "M" for data rows 2-5
$ws.Range("H2:H5").Value = "M"

# Update selection to H6
$ws.Range("H6").Select()
